# Update Switzerland Brokerage & Investment Banking capital-structure database.
# The underlying company rows were recomputed; two company rows (3 and 4)
# also swapped order (Swissquote <-> Compagnie Financiere Tradition).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 ("3" / aggregate-ish row)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 0.0799
$ws.Range("E2").Value = 0.6395
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 149.95
$ws.Range("L2").Value = 0.1135210841093194
$ws.Range("M2").Value = 68.21000000000001
$ws.Range("N2").Value = 0.02820576437993633
$ws.Range("O2").Value = 0.4548849616538847
$ws.Range("P2").Value = 53.40000000000001
$ws.Range("Q2").Value = 0.02208162758962908
$ws.Range("R2").Value = 0.3561187062354119
$ws.Range("S2").Value = 14.81
$ws.Range("T2").Value = 0.2171235889165812
$ws.Range("U2").Value = 3447.05
$ws.Range("V2").Value = 1.425402142000579
$ws.Range("W2").Value = 0.2058212058212058
$ws.Range("X2").Value = 0.03627442344356857
$ws.Range("Y2").Value = 0.1695467823776373
$ws.Range("Z2").Value = -0.4357651374694018
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.0302374292035237
$ws.Range("AC2").Value = -0.0302374292035237
$ws.Range("AD2").Value = 481.3
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 481.3
$ws.Range("AG2").Value = -2965.75
$ws.Range("AH2").Value = 0.1659884121947855
$ws.Range("AI2").Value = 0.3305404848568093
$ws.Range("AJ2").Value = 5.417389715955792
$ws.Range("AK2").Value = 1.489615510183581
$ws.Range("AM2").Value = -6.319999999999999
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# ---------------------------------------------------------------------------
# Row 3 (now Swissquote Group Holding Ltd (SWX:SQN))
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Swissquote Group Holding Ltd (SWX:SQN)"
$ws.Range("D3").Value = 0.133
$ws.Range("E3").Value = 1.101
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 77.09999999999999
$ws.Range("L3").Value = 0.2488702388637831
$ws.Range("M3").Value = 19.46
$ws.Range("N3").Value = 0.01348953278802163
$ws.Range("O3").Value = 0.2523994811932556
$ws.Range("P3").Value = 15.7
$ws.Range("Q3").Value = 0.01088312768612228
$ws.Range("R3").Value = 0.2036316472114137
$ws.Range("S3").Value = 3.760000000000002
$ws.Range("T3").Value = 0.1932168550873588
$ws.Range("U3").Value = 3088.8
$ws.Range("V3").Value = 2.141134063496465
$ws.Range("W3").Value = 0.2108285479901559
$ws.Range("X3").Value = 0.02953243135499795
$ws.Range("Y3").Value = 0.1812961166351579
$ws.Range("Z3").Value = -0.08630488076665924
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.02911652075236308
$ws.Range("AC3").Value = -0.02911652075236308
$ws.Range("AD3").Value = 47
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 47
$ws.Range("AG3").Value = -3041.8
$ws.Range("AH3").Value = 0.03155209452201933
$ws.Range("AI3").Value = 0.09763190693809722
$ws.Range("AJ3").Value = 1.902076038019009
$ws.Range("AK3").Value = 1.166602746030529
$ws.Range("AM3").Value = 0
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AQ3").ClearContents()

# ---------------------------------------------------------------------------
# Row 4 (now Compagnie Financiere Tradition SA (SWX:CFT))
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Compagnie Financière Tradition SA (SWX:CFT)"
$ws.Range("D4").Value = 0.0268
$ws.Range("E4").Value = 0.178
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 79.2
$ws.Range("L4").Value = 0.07893163244967112
$ws.Range("M4").Value = 47.58000000000001
$ws.Range("N4").Value = 0.05053637812002125
$ws.Range("O4").Value = 0.6007575757575758
$ws.Range("P4").Value = 37.7
$ws.Range("Q4").Value = 0.04004248539564525
$ws.Range("R4").Value = 0.476010101010101
$ws.Range("S4").Value = 9.880000000000003
$ws.Range("T4").Value = 0.2076502732240438
$ws.Range("U4").Value = 348.5
$ws.Range("V4").Value = 0.370154009559214
$ws.Range("W4").Value = 0.2058212058212058
$ws.Range("X4").Value = 0.03627442344356857
$ws.Range("Y4").Value = 0.1695467823776373
$ws.Range("Z4").Value = 2.32753421479935
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0302374292035237
$ws.Range("AC4").Value = -0.0302374292035237
$ws.Range("AD4").Value = 409.3
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 409.3
$ws.Range("AG4").Value = 60.80000000000001
$ws.Range("AH4").Value = 0.3030056262955286
$ws.Range("AI4").Value = 0.4926576793452095
$ws.Range("AJ4").Value = 0.06066048089394394
$ws.Range("AK4").Value = 0.1260626166286544
$ws.Range("AM4").Value = -1.26
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("AQ4").Value = -0

# ---------------------------------------------------------------------------
# Row 5 (Valartis Group AG (SWX:VLRT))
# ---------------------------------------------------------------------------
$ws.Range("D5").ClearContents()
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -6.35
$ws.Range("L5").Value = -0.8246753246753246
$ws.Range("M5").Value = 1.17
$ws.Range("N5").Value = 0.03421052631578947
$ws.Range("O5").Value = -0.184251968503937
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 1.17
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 9.75
$ws.Range("V5").Value = 0.2850877192982456
$ws.Range("W5").Value = -0.05857933579335792
$ws.Range("X5").Value = 0.04124120121827186
$ws.Range("Y5").Value = -0.0998205370116298
$ws.Range("Z5").Value = 0.06049654305468258
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03073001827136652
$ws.Range("AC5").Value = -0.03073001827136652
$ws.Range("AD5").Value = 25
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 25
$ws.Range("AG5").Value = 15.25
$ws.Range("AH5").Value = 0.4222972972972973
$ws.Range("AI5").Value = 0.1737317581653926
$ws.Range("AJ5").Value = 0.3083923154701719
$ws.Range("AK5").Value = 0.1136787178531495
$ws.Range("AM5").Value = -5.06
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

Write-Host "Switzerland brokerage & investment banking database updated."
